# This workbook is a weekly price log. A new observation is added at the
# top of the data block (row 4), pushing all the existing data rows down
# by one (old row 4 -> new row 5, ..., old row 117 -> new row 118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right above the current row 4 (first data row after the
# 2 existing "2021-11-22" entries). This shifts rows 4:117 down to 5:118
# and automatically grows the sheet dimension to A1:T118.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 44756
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107002
$ws.Range("J4").Value = "Chirimoya"
$ws.Range("K4").Value = "Cultivar IV Región"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 4000
$ws.Range("O4").Value = 4000
$ws.Range("P4").Value = 4000
$ws.Range("Q4").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R4").Value = "Provincia del Elquí"
$ws.Range("S4").Value = 4000
$ws.Range("T4").Value = 1
